$wb = $excel.ActiveWorkbook

# --- Sheet: Diameter of hair ---
$wsHair = $wb.Worksheets.Item("Diameter of hair")

# Fill column D (rows 9-17) with sequential values 7..15
$wsHair.Range("D9").Value = 7
$wsHair.Range("D10").Value = 8
$wsHair.Range("D11").Value = 9
$wsHair.Range("D12").Value = 10
$wsHair.Range("D13").Value = 11
$wsHair.Range("D14").Value = 12
$wsHair.Range("D15").Value = 13
$wsHair.Range("D16").Value = 14
$wsHair.Range("D17").Value = 15

# Remove the now-unused trailing rows (18 and 19), shrinking the sheet's dimension to D1:F17
$wsHair.Rows("18:19").Delete()

# Activate this worksheet (becomes the selected tab) and set the new selection
$wsHair.Activate()
$wsHair.Range("B6").Select()

# --- Sheet: Wavelegnth with ruler --- (loses tabSelected, keeps its own selection)
$wsWave = $wb.Worksheets.Item("Wavelegnth with ruler")
$wsWave.Range("I3:I23").Select()

# --- Sheet: Refraction os Solid --- (selection only changes)
$wsRefr = $wb.Worksheets.Item("Refraction os Solid")
$wsRefr.Range("G8").Select()

# Re-activate "Diameter of hair" last so it is the workbook's active sheet / tab
$wsHair.Activate()
